# Fix transmission multiplier misspecification:
#  - add a new row 4 "no_vaccination", copied from row 3 (via formulas referencing row 3),
#    except target.coverage (L) = 0 and time.to.coverage (M) = 180 (plain literal values)
#  - rename row 2 ("30p_in_12m_12wk_interval") to "30p_in_6m_12wk_interval"
#    and update its time.to.coverage (M) from 360 to 180
#  - rename row 3 ("30p_in_12m_lumped") to "30p_in_6m_lumped"
#    and update its time.to.coverage (M) from 360 to 180

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build new row 4 by copying row 3 (formats + formulas referencing row 3).
$ws.Range("C3:AM3").Copy()
$ws.Range("C4").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B4").Formula = "=B3"
$ws.Range("C4").Formula = "=C3"
$ws.Range("D4:AM4").Formula = "=D3"

# Override the two cells that differ for the "no vaccination" scenario.
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 180

# Set the run name for the new row first (it claims the freed string slot).
$ws.Range("A4").Value = "no_vaccination"

# Rename the existing scenarios from 12-month to 6-month dosing interval,
# and fix their time.to.coverage (360 -> 180).
$ws.Range("A2").Value = "30p_in_6m_12wk_interval"
$ws.Range("M2").Value = 180

$ws.Range("A3").Value = "30p_in_6m_lumped"
$ws.Range("M3").Value = 180

$ws.Calculate() | Out-Null
